$wb = $excel.ActiveWorkbook

$alc = $wb.Worksheets.Item("ALC")
$arm = $wb.Worksheets.Item("ARM")
$bsm = $wb.Worksheets.Item("BSM")
$crp = $wb.Worksheets.Item("CRP")
$cul = $wb.Worksheets.Item("CUL")
$gsm = $wb.Worksheets.Item("GSM")
$ltw = $wb.Worksheets.Item("LTW")
$wvr = $wb.Worksheets.Item("WVR")

# ---- ALC ----
# row 4 (G4=5470)
$alc.Range("H4").Value = 273.75
$alc.Range("I4").Value = 273.75
$alc.Range("K4").Value = 273.75
$alc.Range("M4").Value = -159.75
# row 19 (G19=7015)
$alc.Range("H19").Value = 1600.2354
$alc.Range("J19").Value = 1922.2222
$alc.Range("L19").Value = 1922.2222
$alc.Range("N19").Value = -2272.2222
# row 43 (G43=5472)
$alc.Range("H43").Value = 1230.4615
$alc.Range("I43").Value = 1228.5714
$alc.Range("J43").Value = 1232.6666
$alc.Range("K43").Value = 1228.5714
$alc.Range("L43").Value = 1232.6666
$alc.Range("M43").Value = -1159.5714
$alc.Range("N43").Value = -1370.6666
# row 51 (G51=5486)
$alc.Range("H51").Value = 5872.909
$alc.Range("I51").Value = 2500
$alc.Range("J51").Value = 6622.4443
$alc.Range("K51").Value = 2500
$alc.Range("L51").Value = 6622.4443
$alc.Range("M51").Value = -2016
$alc.Range("N51").Value = -7590.4443
# row 106 (G106=19903)
$alc.Range("H106").Value = 562443.25
$alc.Range("I106").Value = 5000524
$alc.Range("J106").Value = 7683.125
$alc.Range("K106").Value = 5000524
$alc.Range("L106").Value = 7683.125
$alc.Range("M106").Value = -4999893
$alc.Range("N106").Value = -8945.125
# row 135 (G135=44047)
$alc.Range("H135").Value = 16667164
$alc.Range("I135").Value = 323.4889
$alc.Range("J135").Value = 66667690
$alc.Range("K135").Value = 2911.4001
$alc.Range("L135").Value = 600009210
$alc.Range("M135").Value = -376.4000999999998
$alc.Range("N135").Value = -600014280
# row 137 (G137=44013)
$alc.Range("H137").Value = 1183229.8
$alc.Range("I137").Value = 1201.6285
$alc.Range("J137").Value = 4630812
$alc.Range("K137").Value = 3604.8855
$alc.Range("L137").Value = 13892436
$alc.Range("M137").Value = -1054.8855
$alc.Range("N137").Value = -13897536
# row 138 (G138=44169)
$alc.Range("H138").Value = 2301158.2
$alc.Range("I138").Value = 921.7273
$alc.Range("J138").Value = 4654888.5
$alc.Range("K138").Value = 2765.1819
$alc.Range("L138").Value = 13964665.5
$alc.Range("M138").Value = 2374.8181
$alc.Range("N138").Value = -13974945.5

# ---- ARM ----
# row 32 (G32=44147)
$arm.Range("H32").Value = 1837.89
$arm.Range("I32").Value = 1837.89
$arm.Range("J32").Value = 0
$arm.Range("K32").Value = 1837.89
$arm.Range("L32").Value = 0
$arm.Range("M32").Value = -1550.89
$arm.Range("N32").ClearContents()
# row 52 (G52=27146)
$arm.Range("H52").Value = 24000
$arm.Range("J52").Value = 24000
$arm.Range("L52").Value = 24000
$arm.Range("N52").Value = -24636
# row 61 (G61=43999)
$arm.Range("H61").Value = 1207.1951
$arm.Range("I61").Value = 1063.1228
$arm.Range("J61").Value = 1535.68
$arm.Range("K61").Value = 1063.1228
$arm.Range("L61").Value = 1535.68
$arm.Range("M61").Value = -851.1228000000001
$arm.Range("N61").Value = -1959.68
# row 136 (G136=43999)
$arm.Range("H136").Value = 1207.1951
$arm.Range("I136").Value = 1063.1228
$arm.Range("J136").Value = 1535.68
$arm.Range("K136").Value = 3189.3684
$arm.Range("L136").Value = 4607.04
$arm.Range("M136").Value = -639.3684000000003
$arm.Range("N136").Value = -9707.040000000001

# ---- BSM ----
# row 9 (G9=1648)
$bsm.Range("H9").Value = 27480
$bsm.Range("J9").Value = 27480
$bsm.Range("L9").Value = 27480
$bsm.Range("N9").Value = -27816

# ---- CRP ----
# row 31 (G31=44023)
$crp.Range("H31").Value = 8476022
$crp.Range("I31").Value = 1136.5128
$crp.Range("J31").Value = 25002050
$crp.Range("K31").Value = 1136.5128
$crp.Range("L31").Value = 25002050
$crp.Range("M31").Value = -841.5128
$crp.Range("N31").Value = -25002640
# row 34 (G34=44023)
$crp.Range("H34").Value = 8476022
$crp.Range("I34").Value = 1136.5128
$crp.Range("J34").Value = 25002050
$crp.Range("K34").Value = 1136.5128
$crp.Range("L34").Value = 25002050
$crp.Range("M34").Value = -934.5128
$crp.Range("N34").Value = -25002454
# row 58 (G58=44021)
$crp.Range("H58").Value = 4140.8184
$crp.Range("I58").Value = 4871.0386
$crp.Range("J58").Value = 1428.5714
$crp.Range("K58").Value = 4871.0386
$crp.Range("L58").Value = 1428.5714
$crp.Range("M58").Value = -4668.0386
$crp.Range("N58").Value = -1834.5714
# row 99 (G99=36198)
$crp.Range("H99").Value = 5635.5625
$crp.Range("I99").Value = 5018.1816
$crp.Range("J99").Value = 6993.8
$crp.Range("K99").Value = 5018.1816
$crp.Range("L99").Value = 6993.8
$crp.Range("M99").Value = -3520.1816
$crp.Range("N99").Value = -9989.799999999999
# row 126 (G126=36198)
$crp.Range("H126").Value = 5635.5625
$crp.Range("I126").Value = 5018.1816
$crp.Range("J126").Value = 6993.8
$crp.Range("K126").Value = 15054.5448
$crp.Range("L126").Value = 20981.4
$crp.Range("M126").Value = -12584.5448
$crp.Range("N126").Value = -25921.4
# row 132 (G132=44019)
$crp.Range("H132").Value = 863327.0600000001
$crp.Range("I132").Value = 2001.2051
$crp.Range("J132").Value = 9261254
$crp.Range("K132").Value = 6003.615299999999
$crp.Range("L132").Value = 27783762
$crp.Range("M132").Value = -3473.615299999999
$crp.Range("N132").Value = -27788822
# row 136 (G136=44021)
$crp.Range("H136").Value = 4140.8184
$crp.Range("I136").Value = 4871.0386
$crp.Range("J136").Value = 1428.5714
$crp.Range("K136").Value = 14613.1158
$crp.Range("L136").Value = 4285.7142
$crp.Range("M136").Value = -12063.1158
$crp.Range("N136").Value = -9385.7142

# ---- CUL ----
# row 4 (G4=4650)
$cul.Range("H4").Value = 1015.4375
$cul.Range("I4").Value = 309.4
$cul.Range("J4").Value = 1336.3636
$cul.Range("K4").Value = 928.1999999999999
$cul.Range("L4").Value = 4009.0908
$cul.Range("M4").Value = -816.1999999999999
$cul.Range("N4").Value = -4233.0908
# row 33 (G33=4867)
$cul.Range("H33").Value = 109.5
$cul.Range("I33").Value = 48.5
$cul.Range("J33").Value = 150.16667
$cul.Range("K33").Value = 291
$cul.Range("L33").Value = 901.0000200000001
$cul.Range("M33").Value = -8
$cul.Range("N33").Value = -1467.00002
# row 107 (G107=27838)
$cul.Range("H107").Value = 1398292.2
$cul.Range("I107").Value = 333.33334
$cul.Range("J107").Value = 1548073.6
$cul.Range("K107").Value = 1000.00002
$cul.Range("L107").Value = 4644220.800000001
$cul.Range("M107").Value = 919.9999799999999
$cul.Range("N107").Value = -4648060.800000001
# row 121 (G121=27878)
$cul.Range("H121").Value = 2431314.8
$cul.Range("I121").Value = 347.5
$cul.Range("J121").Value = 3241637.2
$cul.Range("K121").Value = 1042.5
$cul.Range("L121").Value = 9724911.600000001
$cul.Range("M121").Value = 267.5
$cul.Range("N121").Value = -9727531.600000001
# row 122 (G122=36078)
$cul.Range("H122").Value = 645.5
$cul.Range("I122").Value = 544.75
$cul.Range("K122").Value = 4902.75
$cul.Range("M122").Value = -2452.75
# row 131 (G131=36060)
$cul.Range("H131").Value = 930.4
$cul.Range("I131").Value = 330
$cul.Range("J131").Value = 936.46466
$cul.Range("K131").Value = 990
$cul.Range("L131").Value = 2809.39398
$cul.Range("M131").Value = 4050
$cul.Range("N131").Value = -12889.39398
# row 137 (G137=44088)
$cul.Range("H137").Value = 15385976
$cul.Range("I137").Value = 1900.7273
$cul.Range("J137").Value = 20514000
$cul.Range("K137").Value = 5702.1819
$cul.Range("L137").Value = 61542000
$cul.Range("M137").Value = -602.1818999999996
$cul.Range("N137").Value = -61552200

# ---- GSM ----
# row 102 (G102=36169)
$gsm.Range("H102").Value = 1930.8667
$gsm.Range("I102").Value = 1761.25
$gsm.Range("J102").Value = 2609.3333
$gsm.Range("K102").Value = 1761.25
$gsm.Range("L102").Value = 2609.3333
$gsm.Range("M102").Value = -139.25
$gsm.Range("N102").Value = -5853.3333
# row 121 (G121=26338)
$gsm.Range("H121").Value = 20105.666
$gsm.Range("J121").Value = 20105.666
$gsm.Range("L121").Value = 20105.666
$gsm.Range("N121").Value = -23599.666
# row 126 (G126=36184)
$gsm.Range("H126").Value = 3752.2632
$gsm.Range("I126").Value = 3956.6924
$gsm.Range("J126").Value = 3309.3333
$gsm.Range("K126").Value = 11870.0772
$gsm.Range("L126").Value = 9927.999899999999
$gsm.Range("M126").Value = -9400.0772
$gsm.Range("N126").Value = -14867.9999

# ---- LTW ----
# row 122 (G122=36247)
$ltw.Range("H122").Value = 8750
$ltw.Range("I122").Value = 20000
$ltw.Range("K122").Value = 60000
$ltw.Range("M122").Value = -57550
# row 132 (G132=44058)
$ltw.Range("H132").Value = 3422.127
$ltw.Range("I132").Value = 3249.1333
$ltw.Range("J132").Value = 3854.611
$ltw.Range("K132").Value = 9747.3999
$ltw.Range("L132").Value = 11563.833
$ltw.Range("M132").Value = -7217.3999
$ltw.Range("N132").Value = -16623.833
# row 136 (G136=44060)
$ltw.Range("H136").Value = 1443.9508
$ltw.Range("I136").Value = 849.3913
$ltw.Range("K136").Value = 2548.1739
$ltw.Range("M136").Value = 1.826100000000224

# ---- WVR ----
# row 39 (G39=3106)
$wvr.Range("H39").Value = 16666.666
$wvr.Range("J39").Value = 16666.666
$wvr.Range("L39").Value = 16666.666
$wvr.Range("N39").Value = -17492.666
# row 42 (G42=3372)
$wvr.Range("H42").Value = 22105
$wvr.Range("J42").Value = 22105
$wvr.Range("L42").Value = 22105
$wvr.Range("N42").Value = -22861
# row 112 (G112=25836)
$wvr.Range("H112").Value = 30000
$wvr.Range("J112").Value = 30000
$wvr.Range("L112").Value = 30000
$wvr.Range("N112").Value = -32954
# row 122 (G122=36208)
$wvr.Range("H122").Value = 323252.25
$wvr.Range("I122").Value = 90004
$wvr.Range("J122").Value = 401001.66
$wvr.Range("K122").Value = 270012
$wvr.Range("L122").Value = 1203004.98
$wvr.Range("M122").Value = -267562
$wvr.Range("N122").Value = -1207904.98
# row 132 (G132=44029)
$wvr.Range("H132").Value = 2711.0227
$wvr.Range("I132").Value = 3002.8333
$wvr.Range("J132").Value = 2085.7144
$wvr.Range("K132").Value = 9008.499899999999
$wvr.Range("L132").Value = 6257.1432
$wvr.Range("M132").Value = -6478.499899999999
$wvr.Range("N132").Value = -11317.1432
# row 136 (G136=44031)
$wvr.Range("H136").Value = 1997.5646
$wvr.Range("I136").Value = 1619.1837
$wvr.Range("J136").Value = 3423.7693
$wvr.Range("K136").Value = 4857.551100000001
$wvr.Range("L136").Value = 10271.3079
$wvr.Range("M136").Value = -2307.551100000001
$wvr.Range("N136").Value = -15371.3079

